$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save current (pre-edit) values of row 4, 5, 6, 7 for the columns that change
$row4 = @{
    D = $ws.Range("D4").Value2
    M = $ws.Range("M4").Value2
    N = $ws.Range("N4").Value2
    O = $ws.Range("O4").Value2
    P = $ws.Range("P4").Value2
    S = $ws.Range("S4").Value2
}
$row5 = @{
    D = $ws.Range("D5").Value2
    M = $ws.Range("M5").Value2
    N = $ws.Range("N5").Value2
    O = $ws.Range("O5").Value2
    P = $ws.Range("P5").Value2
    S = $ws.Range("S5").Value2
}
$row6 = @{
    D = $ws.Range("D6").Value2
    M = $ws.Range("M6").Value2
    N = $ws.Range("N6").Value2
    O = $ws.Range("O6").Value2
    P = $ws.Range("P6").Value2
    S = $ws.Range("S6").Value2
}
$row7 = @{
    D = $ws.Range("D7").Value2
    M = $ws.Range("M7").Value2
    N = $ws.Range("N7").Value2
    O = $ws.Range("O7").Value2
    P = $ws.Range("P7").Value2
    S = $ws.Range("S7").Value2
}

# Swap row 4 <-> row 6
$ws.Range("D4").Value = $row6.D
$ws.Range("M4").Value = $row6.M
$ws.Range("N4").Value = $row6.N
$ws.Range("O4").Value = $row6.O
$ws.Range("P4").Value = $row6.P
$ws.Range("S4").Value = $row6.S

$ws.Range("D6").Value = $row4.D
$ws.Range("M6").Value = $row4.M
$ws.Range("N6").Value = $row4.N
$ws.Range("O6").Value = $row4.O
$ws.Range("P6").Value = $row4.P
$ws.Range("S6").Value = $row4.S

# Swap row 5 <-> row 7
$ws.Range("D5").Value = $row7.D
$ws.Range("M5").Value = $row7.M
$ws.Range("N5").Value = $row7.N
$ws.Range("O5").Value = $row7.O
$ws.Range("P5").Value = $row7.P
$ws.Range("S5").Value = $row7.S

$ws.Range("D7").Value = $row5.D
$ws.Range("M7").Value = $row5.M
$ws.Range("N7").Value = $row5.N
$ws.Range("O7").Value = $row5.O
$ws.Range("P7").Value = $row5.P
$ws.Range("S7").Value = $row5.S
